# Applies the Linea 141 schedule refresh (scrape @ 05:52:07) described in the commit
# "Horarios actualizados Linea 141 - 821" to all three worksheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = "Última actualización: 05:52:07"
$ws.Cells.Item(3, 1).Value = "Total filas: 54"

$ws.Cells.Item(29, 1).Value = "05:52:07"
$ws.Cells.Item(29, 2).Value = "05:54"
$ws.Cells.Item(29, 3).Value = "10_OLMOS"
$ws.Cells.Item(29, 4).Value = 2
$ws.Cells.Item(29, 5).Value = "LP1912"
$ws.Cells.Item(32, 1).Value = "05:52:07"
$ws.Cells.Item(32, 2).Value = "06:12"
$ws.Cells.Item(32, 3).Value = "215A_EL PATO"
$ws.Cells.Item(32, 4).Value = 20
$ws.Cells.Item(32, 5).Value = "LP1912"
$ws.Cells.Item(33, 1).Value = "04:48:57"
$ws.Cells.Item(33, 2).Value = "06:13"
$ws.Cells.Item(33, 3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(33, 4).Value = 85
$ws.Cells.Item(33, 5).Value = "LP1912"
$ws.Cells.Item(34, 1).Value = "05:52:07"
$ws.Cells.Item(34, 2).Value = "06:14"
$ws.Cells.Item(34, 3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(34, 4).Value = 22
$ws.Cells.Item(34, 5).Value = "LP1912"
$ws.Cells.Item(35, 1).Value = "04:48:57"
$ws.Cells.Item(35, 2).Value = "06:20"
$ws.Cells.Item(35, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(35, 4).Value = 92
$ws.Cells.Item(35, 5).Value = "LP1912"
$ws.Cells.Item(36, 1).Value = "05:52:07"
$ws.Cells.Item(36, 2).Value = "06:21"
$ws.Cells.Item(36, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(36, 4).Value = 29
$ws.Cells.Item(36, 5).Value = "LP1912"
$ws.Cells.Item(37, 1).Value = "04:48:57"
$ws.Cells.Item(37, 2).Value = "06:26"
$ws.Cells.Item(37, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(37, 4).Value = 98
$ws.Cells.Item(37, 5).Value = "LP1912"
$ws.Cells.Item(38, 1).Value = "05:52:07"
$ws.Cells.Item(38, 2).Value = "06:27"
$ws.Cells.Item(38, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(38, 4).Value = 35
$ws.Cells.Item(38, 5).Value = "LP1912"
$ws.Cells.Item(39, 1).Value = "04:48:57"
$ws.Cells.Item(39, 2).Value = "06:29"
$ws.Cells.Item(39, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(39, 4).Value = 101
$ws.Cells.Item(39, 5).Value = "LP1912"
$ws.Cells.Item(40, 1).Value = "05:52:07"
$ws.Cells.Item(40, 2).Value = "06:30"
$ws.Cells.Item(40, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(40, 4).Value = 38
$ws.Cells.Item(40, 5).Value = "LP1912"
$ws.Cells.Item(41, 1).Value = "05:52:07"
$ws.Cells.Item(41, 2).Value = "06:31"
$ws.Cells.Item(41, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(41, 4).Value = 39
$ws.Cells.Item(41, 5).Value = "LP1912"
$ws.Cells.Item(42, 1).Value = "04:48:57"
$ws.Cells.Item(42, 2).Value = "06:43"
$ws.Cells.Item(42, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(42, 4).Value = 115
$ws.Cells.Item(42, 5).Value = "LP1912"
$ws.Cells.Item(43, 1).Value = "05:52:07"
$ws.Cells.Item(43, 2).Value = "06:44"
$ws.Cells.Item(43, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(43, 4).Value = 52
$ws.Cells.Item(43, 5).Value = "LP1912"
$ws.Cells.Item(44, 1).Value = "05:21:16"
$ws.Cells.Item(44, 2).Value = "06:46"
$ws.Cells.Item(44, 3).Value = "215C_EL PATO"
$ws.Cells.Item(44, 4).Value = 85
$ws.Cells.Item(44, 5).Value = "LP1912"
$ws.Cells.Item(45, 1).Value = "05:52:07"
$ws.Cells.Item(45, 2).Value = "06:47"
$ws.Cells.Item(45, 3).Value = "215C_EL PATO"
$ws.Cells.Item(45, 4).Value = 55
$ws.Cells.Item(45, 5).Value = "LP1912"
$ws.Cells.Item(46, 1).Value = "05:52:07"
$ws.Cells.Item(46, 2).Value = "07:00"
$ws.Cells.Item(46, 3).Value = "14_ABASTO"
$ws.Cells.Item(46, 4).Value = 68
$ws.Cells.Item(46, 5).Value = "LP1912"
$ws.Cells.Item(47, 1).Value = "05:52:07"
$ws.Cells.Item(47, 2).Value = "07:05"
$ws.Cells.Item(47, 3).Value = "15_ABASTO"
$ws.Cells.Item(47, 4).Value = 73
$ws.Cells.Item(47, 5).Value = "LP1912"
$ws.Cells.Item(48, 1).Value = "05:52:07"
$ws.Cells.Item(48, 2).Value = "07:05"
$ws.Cells.Item(48, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(48, 4).Value = 73
$ws.Cells.Item(48, 5).Value = "LP1912"
$ws.Cells.Item(49, 1).Value = "05:52:07"
$ws.Cells.Item(49, 2).Value = "07:07"
$ws.Cells.Item(49, 3).Value = "225_GOMEZ"
$ws.Cells.Item(49, 4).Value = 75
$ws.Cells.Item(49, 5).Value = "LP1912"
$ws.Cells.Item(50, 1).Value = "05:21:16"
$ws.Cells.Item(50, 2).Value = "07:11"
$ws.Cells.Item(50, 3).Value = "215A_EL PATO"
$ws.Cells.Item(50, 4).Value = 110
$ws.Cells.Item(50, 5).Value = "LP1912"
$ws.Cells.Item(51, 1).Value = "05:52:07"
$ws.Cells.Item(51, 2).Value = "07:12"
$ws.Cells.Item(51, 3).Value = "215A_EL PATO"
$ws.Cells.Item(51, 4).Value = 80
$ws.Cells.Item(51, 5).Value = "LP1912"
$ws.Cells.Item(52, 1).Value = "05:52:07"
$ws.Cells.Item(52, 2).Value = "07:16"
$ws.Cells.Item(52, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(52, 4).Value = 84
$ws.Cells.Item(52, 5).Value = "LP1912"
$ws.Cells.Item(53, 1).Value = "05:52:07"
$ws.Cells.Item(53, 2).Value = "07:21"
$ws.Cells.Item(53, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(53, 4).Value = 89
$ws.Cells.Item(53, 5).Value = "LP1912"
$ws.Cells.Item(54, 1).Value = "05:52:07"
$ws.Cells.Item(54, 2).Value = "07:23"
$ws.Cells.Item(54, 3).Value = "10_OLMOS"
$ws.Cells.Item(54, 4).Value = 91
$ws.Cells.Item(54, 5).Value = "LP1912"
$ws.Cells.Item(55, 1).Value = "05:52:07"
$ws.Cells.Item(55, 2).Value = "07:32"
$ws.Cells.Item(55, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(55, 4).Value = 100
$ws.Cells.Item(55, 5).Value = "LP1912"
$ws.Cells.Item(56, 1).Value = "05:52:07"
$ws.Cells.Item(56, 2).Value = "07:32"
$ws.Cells.Item(56, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(56, 4).Value = 100
$ws.Cells.Item(56, 5).Value = "LP1912"
$ws.Cells.Item(57, 1).Value = "05:52:07"
$ws.Cells.Item(57, 2).Value = "07:32"
$ws.Cells.Item(57, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(57, 4).Value = 100
$ws.Cells.Item(57, 5).Value = "LP1912"
$ws.Cells.Item(58, 1).Value = "05:52:07"
$ws.Cells.Item(58, 2).Value = "07:37"
$ws.Cells.Item(58, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(58, 4).Value = 105
$ws.Cells.Item(58, 5).Value = "LP1912"
$ws.Cells.Item(59, 1).Value = "05:52:07"
$ws.Cells.Item(59, 2).Value = "07:48"
$ws.Cells.Item(59, 3).Value = "14_ABASTO"
$ws.Cells.Item(59, 4).Value = 116
$ws.Cells.Item(59, 5).Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(2, 1).Value = "Última actualización: 05:52:07"
$ws.Cells.Item(3, 1).Value = "Total filas: 14"

$ws.Cells.Item(15, 1).Value = "05:52:07"
$ws.Cells.Item(15, 2).Value = "06:12"
$ws.Cells.Item(15, 3).Value = "215A_EL PATO"
$ws.Cells.Item(15, 4).Value = 20
$ws.Cells.Item(15, 5).Value = "LP1912"
$ws.Cells.Item(16, 1).Value = "05:21:16"
$ws.Cells.Item(16, 2).Value = "06:46"
$ws.Cells.Item(16, 3).Value = "215C_EL PATO"
$ws.Cells.Item(16, 4).Value = 85
$ws.Cells.Item(16, 5).Value = "LP1912"
$ws.Cells.Item(17, 1).Value = "05:52:07"
$ws.Cells.Item(17, 2).Value = "06:47"
$ws.Cells.Item(17, 3).Value = "215C_EL PATO"
$ws.Cells.Item(17, 4).Value = 55
$ws.Cells.Item(17, 5).Value = "LP1912"
$ws.Cells.Item(18, 1).Value = "05:21:16"
$ws.Cells.Item(18, 2).Value = "07:11"
$ws.Cells.Item(18, 3).Value = "215A_EL PATO"
$ws.Cells.Item(18, 4).Value = 110
$ws.Cells.Item(18, 5).Value = "LP1912"
$ws.Cells.Item(19, 1).Value = "05:52:07"
$ws.Cells.Item(19, 2).Value = "07:12"
$ws.Cells.Item(19, 3).Value = "215A_EL PATO"
$ws.Cells.Item(19, 4).Value = 80
$ws.Cells.Item(19, 5).Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)

$ws.Cells.Item(2, 1).Value = "Última actualización: 05:52:07"
$ws.Cells.Item(3, 1).Value = "Total filas: 10"

$ws.Cells.Item(11, 1).Value = "05:52:07"
$ws.Cells.Item(11, 2).Value = "06:13"
$ws.Cells.Item(11, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(11, 4).Value = 21
$ws.Cells.Item(11, 5).Value = "L6173"
$ws.Cells.Item(12, 1).Value = "04:48:57"
$ws.Cells.Item(12, 2).Value = "06:32"
$ws.Cells.Item(12, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(12, 4).Value = 104
$ws.Cells.Item(12, 5).Value = "L6203"
$ws.Cells.Item(13, 1).Value = "05:52:07"
$ws.Cells.Item(13, 2).Value = "06:33"
$ws.Cells.Item(13, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(13, 4).Value = 41
$ws.Cells.Item(13, 5).Value = "L6203"
$ws.Cells.Item(14, 1).Value = "05:52:07"
$ws.Cells.Item(14, 2).Value = "07:00"
$ws.Cells.Item(14, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(14, 4).Value = 68
$ws.Cells.Item(14, 5).Value = "L6173"
$ws.Cells.Item(15, 1).Value = "05:52:07"
$ws.Cells.Item(15, 2).Value = "07:35"
$ws.Cells.Item(15, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(15, 4).Value = 103
$ws.Cells.Item(15, 5).Value = "L6173"

